$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Interventions target population")

# Normalize B2's number format so it shares the same base format (General/"0")
# as the rest of the row before centering the whole block - this keeps every
# cell in B2:G8 collapsing onto a single shared style once alignment is applied.
$ws.Range("B2").NumberFormat = "general"

$values = @(
    @(0, 0, 0, 1, 1, 0),
    @(0, 0, 1, 1, 1, 0),
    @(0, 1, 1, 1, 1, 0),
    @(1, 1, 1, 1, 0, 0),
    @(0, 0, 0, 0, 0, 1),
    @(0, 0, 0, 0, 0, 1),
    @(0, 0, 0, 0, 0, 1)
)

for ($i = 0; $i -lt 7; $i++) {
    $row = 2 + $i
    for ($col = 2; $col -le 7; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$i][$col - 2]
    }
}

# Center-align the whole populated target-population matrix.
$ws.Range("B2:G8").HorizontalAlignment = -4108
